# Add a new "UptimeValue" column (X) to the datasheet, to the right of the
# existing "DownTime" column (W).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new column (X / index 24) its own width, matching the source edit.
$ws.Columns.Item(24).ColumnWidth = 24.5

# X1: header cell - reuse the header formatting used by the other header
# cells in row 1 (same style as W1), then write the header text.
$ws.Cells.Item(1, 23).Copy()
$ws.Cells.Item(1, 24).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 24).Value = "UptimeValue"

# X2: data cell - reuse the formatting already used for the "GuranteedUptime"
# value cell (T1, Monaco font style), then write the numeric value.
$ws.Cells.Item(1, 20).Copy()
$ws.Cells.Item(2, 24).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(2, 24).Value = 98.33

$excel.CutCopyMode = $false

# Move/leave the selection where the author left it after the edit.
$ws.Range("W10").Select()
